$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 262 (data index 260) ---
$ws.Range("D262:G262").NumberFormat = "@"
$ws.Range("I262").NumberFormat = "@"
$ws.Range("K262:L262").NumberFormat = "@"
$ws.Cells.Item(262,4).Value = "25047.56000000"
$ws.Cells.Item(262,5).Value = "22664.69000000"
$ws.Cells.Item(262,6).Value = "24305.24000000"
$ws.Cells.Item(262,7).Value = "1251083.26468000"
$ws.Cells.Item(262,9).Value = "30077619286.84656770"
$ws.Cells.Item(262,10).Value = 43311786
$ws.Cells.Item(262,11).Value = "625522.69768000"
$ws.Cells.Item(262,12).Value = "15039436995.29061950"
$ws.Range("D262:G262").ClearFormats()
$ws.Range("I262").ClearFormats()
$ws.Range("K262:L262").ClearFormats()

# --- Append new rows 263-271 ---
$ws.Range("A262").Copy($ws.Range("A263:A271"))

$ws.Cells.Item(263,1).Value = 261
$ws.Cells.Item(263,2).Value = 1660521600000
$ws.Range("C263:G263").NumberFormat = "@"
$ws.Cells.Item(263,3).Value = "24305.25000000"
$ws.Cells.Item(263,4).Value = "25211.32000000"
$ws.Cells.Item(263,5).Value = "20761.90000000"
$ws.Cells.Item(263,6).Value = "21515.61000000"
$ws.Cells.Item(263,7).Value = "1402957.39876000"
$ws.Cells.Item(263,8).Value = 1661126399999
$ws.Range("I263").NumberFormat = "@"
$ws.Cells.Item(263,9).Value = "32016009263.43493180"
$ws.Cells.Item(263,10).Value = 46017032
$ws.Range("K263:M263").NumberFormat = "@"
$ws.Cells.Item(263,11).Value = "698823.59671000"
$ws.Cells.Item(263,12).Value = "15948153490.60180970"
$ws.Cells.Item(263,13).Value = "0"
$ws.Range("C263:G263").ClearFormats()
$ws.Range("I263").ClearFormats()
$ws.Range("K263:M263").ClearFormats()

$ws.Cells.Item(264,1).Value = 262
$ws.Cells.Item(264,2).Value = 1661126400000
$ws.Range("C264:G264").NumberFormat = "@"
$ws.Cells.Item(264,3).Value = "21516.70000000"
$ws.Cells.Item(264,4).Value = "21900.00000000"
$ws.Cells.Item(264,5).Value = "19520.00000000"
$ws.Cells.Item(264,6).Value = "19555.61000000"
$ws.Cells.Item(264,7).Value = "1343190.86000000"
$ws.Cells.Item(264,8).Value = 1661731199999
$ws.Range("I264").NumberFormat = "@"
$ws.Cells.Item(264,9).Value = "28273617081.52032370"
$ws.Cells.Item(264,10).Value = 41486274
$ws.Range("K264:M264").NumberFormat = "@"
$ws.Cells.Item(264,11).Value = "669270.17315000"
$ws.Cells.Item(264,12).Value = "14089524497.01568180"
$ws.Cells.Item(264,13).Value = "0"
$ws.Range("C264:G264").ClearFormats()
$ws.Range("I264").ClearFormats()
$ws.Range("K264:M264").ClearFormats()

$ws.Cells.Item(265,1).Value = 263
$ws.Cells.Item(265,2).Value = 1661731200000
$ws.Range("C265:G265").NumberFormat = "@"
$ws.Cells.Item(265,3).Value = "19555.61000000"
$ws.Cells.Item(265,4).Value = "20576.25000000"
$ws.Cells.Item(265,5).Value = "19540.00000000"
$ws.Cells.Item(265,6).Value = "20000.30000000"
$ws.Cells.Item(265,7).Value = "1527594.84529000"
$ws.Cells.Item(265,8).Value = 1662335999999
$ws.Range("I265").NumberFormat = "@"
$ws.Cells.Item(265,9).Value = "30597230623.48951400"
$ws.Cells.Item(265,10).Value = 38080138
$ws.Range("K265:M265").NumberFormat = "@"
$ws.Cells.Item(265,11).Value = "762324.89492000"
$ws.Cells.Item(265,12).Value = "15270693969.60497590"
$ws.Cells.Item(265,13).Value = "0"
$ws.Range("C265:G265").ClearFormats()
$ws.Range("I265").ClearFormats()
$ws.Range("K265:M265").ClearFormats()

$ws.Cells.Item(266,1).Value = 264
$ws.Cells.Item(266,2).Value = 1662336000000
$ws.Range("C266:G266").NumberFormat = "@"
$ws.Cells.Item(266,3).Value = "20000.30000000"
$ws.Cells.Item(266,4).Value = "21860.00000000"
$ws.Cells.Item(266,5).Value = "18510.77000000"
$ws.Cells.Item(266,6).Value = "21826.87000000"
$ws.Cells.Item(266,7).Value = "2146685.76233000"
$ws.Cells.Item(266,8).Value = 1662940799999
$ws.Range("I266").NumberFormat = "@"
$ws.Cells.Item(266,9).Value = "43460053550.45976670"
$ws.Cells.Item(266,10).Value = 41587411
$ws.Range("K266:M266").NumberFormat = "@"
$ws.Cells.Item(266,11).Value = "1074020.48582000"
$ws.Cells.Item(266,12).Value = "21751877355.13041250"
$ws.Cells.Item(266,13).Value = "0"
$ws.Range("C266:G266").ClearFormats()
$ws.Range("I266").ClearFormats()
$ws.Range("K266:M266").ClearFormats()

$ws.Cells.Item(267,1).Value = 265
$ws.Cells.Item(267,2).Value = 1662940800000
$ws.Range("C267:G267").NumberFormat = "@"
$ws.Cells.Item(267,3).Value = "21826.87000000"
$ws.Cells.Item(267,4).Value = "22799.00000000"
$ws.Cells.Item(267,5).Value = "19320.01000000"
$ws.Cells.Item(267,6).Value = "19416.18000000"
$ws.Cells.Item(267,7).Value = "2218565.59694000"
$ws.Cells.Item(267,8).Value = 1663545599999
$ws.Range("I267").NumberFormat = "@"
$ws.Cells.Item(267,9).Value = "45784151832.43753880"
$ws.Cells.Item(267,10).Value = 41728131
$ws.Range("K267:M267").NumberFormat = "@"
$ws.Cells.Item(267,11).Value = "1107144.73306000"
$ws.Cells.Item(267,12).Value = "22859382949.06557150"
$ws.Cells.Item(267,13).Value = "0"
$ws.Range("C267:G267").ClearFormats()
$ws.Range("I267").ClearFormats()
$ws.Range("K267:M267").ClearFormats()

$ws.Cells.Item(268,1).Value = 266
$ws.Cells.Item(268,2).Value = 1663545600000
$ws.Range("C268:G268").NumberFormat = "@"
$ws.Cells.Item(268,3).Value = "19417.45000000"
$ws.Cells.Item(268,4).Value = "19956.00000000"
$ws.Cells.Item(268,5).Value = "18125.98000000"
$ws.Cells.Item(268,6).Value = "18807.38000000"
$ws.Cells.Item(268,7).Value = "2285541.48793000"
$ws.Cells.Item(268,8).Value = 1664150399999
$ws.Range("I268").NumberFormat = "@"
$ws.Cells.Item(268,9).Value = "43488049829.07041010"
$ws.Cells.Item(268,10).Value = 39408640
$ws.Range("K268:M268").NumberFormat = "@"
$ws.Cells.Item(268,11).Value = "1141577.20385000"
$ws.Cells.Item(268,12).Value = "21724013825.41485850"
$ws.Cells.Item(268,13).Value = "0"
$ws.Range("C268:G268").ClearFormats()
$ws.Range("I268").ClearFormats()
$ws.Range("K268:M268").ClearFormats()

$ws.Cells.Item(269,1).Value = 267
$ws.Cells.Item(269,2).Value = 1664150400000
$ws.Range("C269:G269").NumberFormat = "@"
$ws.Cells.Item(269,3).Value = "18809.13000000"
$ws.Cells.Item(269,4).Value = "20385.86000000"
$ws.Cells.Item(269,5).Value = "18471.28000000"
$ws.Cells.Item(269,6).Value = "19056.80000000"
$ws.Cells.Item(269,7).Value = "2777070.91238000"
$ws.Cells.Item(269,8).Value = 1664755199999
$ws.Range("I269").NumberFormat = "@"
$ws.Cells.Item(269,9).Value = "53761170640.89073340"
$ws.Cells.Item(269,10).Value = 39023576
$ws.Range("K269:M269").NumberFormat = "@"
$ws.Cells.Item(269,11).Value = "1387219.15628000"
$ws.Cells.Item(269,12).Value = "26857951663.90879500"
$ws.Cells.Item(269,13).Value = "0"
$ws.Range("C269:G269").ClearFormats()
$ws.Range("I269").ClearFormats()
$ws.Range("K269:M269").ClearFormats()

$ws.Cells.Item(270,1).Value = 268
$ws.Cells.Item(270,2).Value = 1664755200000
$ws.Range("C270:G270").NumberFormat = "@"
$ws.Cells.Item(270,3).Value = "19057.74000000"
$ws.Cells.Item(270,4).Value = "20475.00000000"
$ws.Cells.Item(270,5).Value = "18959.68000000"
$ws.Cells.Item(270,6).Value = "19439.02000000"
$ws.Cells.Item(270,7).Value = "1690215.44019000"
$ws.Cells.Item(270,8).Value = 1665359999999
$ws.Range("I270").NumberFormat = "@"
$ws.Cells.Item(270,9).Value = "33498469288.87016030"
$ws.Cells.Item(270,10).Value = 28764711
$ws.Range("K270:M270").NumberFormat = "@"
$ws.Cells.Item(270,11).Value = "847579.42346000"
$ws.Cells.Item(270,12).Value = "16800184632.09654350"
$ws.Cells.Item(270,13).Value = "0"
$ws.Range("C270:G270").ClearFormats()
$ws.Range("I270").ClearFormats()
$ws.Range("K270:M270").ClearFormats()

$ws.Cells.Item(271,1).Value = 269
$ws.Cells.Item(271,2).Value = 1665360000000
$ws.Range("C271:G271").NumberFormat = "@"
$ws.Cells.Item(271,3).Value = "19439.96000000"
$ws.Cells.Item(271,4).Value = "19951.87000000"
$ws.Cells.Item(271,5).Value = "18190.00000000"
$ws.Cells.Item(271,6).Value = "19175.86000000"
$ws.Cells.Item(271,7).Value = "1404814.74082000"
$ws.Cells.Item(271,8).Value = 1665964799999
$ws.Range("I271").NumberFormat = "@"
$ws.Cells.Item(271,9).Value = "26949592667.33566080"
$ws.Cells.Item(271,10).Value = 27587086
$ws.Range("K271:M271").NumberFormat = "@"
$ws.Cells.Item(271,11).Value = "700157.27233000"
$ws.Cells.Item(271,12).Value = "13433257689.75192720"
$ws.Cells.Item(271,13).Value = "0"
$ws.Range("C271:G271").ClearFormats()
$ws.Range("I271").ClearFormats()
$ws.Range("K271:M271").ClearFormats()
